# Update the "想去人数" (F column) counts on both the "展览" and "全部类型"
# sheets, which contain duplicate data tables.

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    2  = 278
    3  = 1413
    4  = 168
    9  = 192
    10 = 139
    11 = 4719
    12 = 6980
    16 = 579
    17 = 57
    19 = 1037
    21 = 71
    22 = 2750
    25 = 177
    27 = 382
    28 = 406
    31 = 1649
    32 = 1053
    34 = 655
    37 = 9
    41 = 226
    42 = 654
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
